# Apply updated cryptocurrency price/volume snapshot to the "cryptos" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 32 and 33 swap position (Stellar <-> Filecoin) and get new data ---
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.682"
$ws.Range("E32").Value = "  +5.47%  "

$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09188"
$ws.Range("E33").Value = "  +1.43%  "

# --- Price / Volume(1h) updates for all other rows ---
# NumberFormat is forced to Text ("@") before assigning D-column values so
# that numeric-looking strings (e.g. "328.59") are stored as text, matching
# the original inline-string cell type, instead of being auto-converted to
# numbers by Excel's type inference.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.960.98"
$ws.Range("E2").Value = "  +1.67%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.769.78"
$ws.Range("E3").Value = "  +1.65%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.59"
$ws.Range("E5").Value = "  +1.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9983"
$ws.Range("E6").Value = "  -0.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4480"
$ws.Range("E7").Value = "  -0.97%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3553"
$ws.Range("E8").Value = "  +0.72%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07434"
$ws.Range("E9").Value = "  +0.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.86"
$ws.Range("E10").Value = "  +1.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.107"
$ws.Range("E11").Value = "  +2.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9980"
$ws.Range("E12").Value = "  -0.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.98"
$ws.Range("E13").Value = "  +2.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.040"
$ws.Range("E14").Value = "  +2.28%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.245"
$ws.Range("E15").Value = "  +2.86%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.769.88"
$ws.Range("E16").Value = "  +1.86%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.40"
$ws.Range("E17").Value = "  +2.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001064"
$ws.Range("E18").Value = "  +1.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06432"
$ws.Range("E19").Value = "  +1.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9991"
$ws.Range("E20").Value = "  -0.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.12"
$ws.Range("E21").Value = "  +2.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.775"
$ws.Range("E22").Value = "  +0.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.992.43"
$ws.Range("E23").Value = "  +1.64%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.30"
$ws.Range("E24").Value = "  +1.61%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.108"
$ws.Range("E25").Value = "  +1.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.29"
$ws.Range("E26").Value = "  -0.67%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.39"
$ws.Range("E27").Value = "  +1.90%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.965.21"
$ws.Range("E28").Value = "  +1.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.162"
# E29 is unchanged ("  +6.00%  ")

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.49"
$ws.Range("E30").Value = "  -0.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.115"
$ws.Range("E31").Value = "  +6.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.683"
$ws.Range("E34").Value = "  +1.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.88"
$ws.Range("E35").Value = "  +2.28%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06215"
$ws.Range("E36").Value = "  +4.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02290"
$ws.Range("E37").Value = "  +0.92%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2110"
$ws.Range("E38").Value = "  +2.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6329"
$ws.Range("E39").Value = "  +1.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.969"
$ws.Range("E40").Value = "  +1.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.188"
$ws.Range("E41").Value = "  +0.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.394"
$ws.Range("E42").Value = "  +1.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.890"
$ws.Range("E43").Value = "  +2.49%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.35"
$ws.Range("E44").Value = "  +1.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.745"
$ws.Range("E45").Value = "  +1.24%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5889"
$ws.Range("E46").Value = "  +1.74%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.58"
$ws.Range("E47").Value = "  +0.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.961"
$ws.Range("E48").Value = "  +1.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.139"
$ws.Range("E49").Value = "  +2.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06882"
$ws.Range("E50").Value = "  +0.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.89"
$ws.Range("E51").Value = "  +2.52%  "
